# "Removing comment anonymity - closes #2"
# Insert a new "User" column into the Notes sheet (between Party and
# Question ID) and make the Notes sheet the active/selected tab instead
# of the Introduction sheet.

$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")

# Insert a new column before column C (Question ID moves from C to D, etc.)
$notes.Columns.Item(3).Insert()

# Give the new header its label
$notes.Cells.Item(1, 3).Value2 = "User"

# Put the active selection on the newly inserted column's second row
$notes.Range("C2").Select()

# Make the Notes sheet the active tab (was Introduction before)
$notes.Activate()
